# Merkertest gedaan. Dut't nie
# Update the "Alias memorybit" column (renamed to "Alias memorybits") so
# memory-bit markers that used to be single addresses become 3-bit ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header / table column: "Alias memorybit" -> "Alias memorybits"
$ws.Range("H6").Value = "Alias memorybits"

$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListColumns.Item("Alias memorybit").Name = "Alias memorybits"

# Update the memory-bit marker values for the output rows (UITGANGEN block)
$ws.Range("H23").Value = "%M71-%M73"
$ws.Range("H24").Value = "%M74-%M76"
$ws.Range("H25").Value = "%M60-%M62"
$ws.Range("H26").Value = "%M65-%M67"
$ws.Range("H27").Value = "%M68-%M70"
$ws.Range("H28").Value = "%M77-%M79"
$ws.Range("H29").Value = "%M80-%M83"

# Row 30's memorybit marker is no longer used - clear it
$ws.Range("H30").ClearContents()

# Flag the two unused transport-lane spare outputs as reserve
$ws.Range("I31").Value = "< transportbaan reserve"
$ws.Range("I32").Value = "< transportbaan reserve"

# Column width tweaks (Alias / Alias memorybits columns widened slightly)
$ws.Columns.Item(7).ColumnWidth = 29.75
$ws.Columns.Item(8).ColumnWidth = 17.75

# Scroll/selection state as left by the author after the edit
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("H21").Select()
